$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "Naked Pairs"
$ws.Range("F6").Value = "x"

$ws.Range("C7").Value = "Hidden Pairs"
$ws.Range("D7").Value = "x"

$ws.Range("H6").Select()
